# Insert a new data row at row 1193 (pushing existing rows 1193-1279 down to
# 1194-1280) and populate it with the new "Femacal de La Calera - Papa"
# observation, matching the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1193).Insert()

$ws.Cells.Item(1193, 1).Value  = 3
$ws.Cells.Item(1193, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(1193, 3).Value  = "Coquimbo"
$ws.Cells.Item(1193, 4).Value  = 45265
$ws.Cells.Item(1193, 5).Value  = 5
$ws.Cells.Item(1193, 6).Value  = 100114001
$ws.Cells.Item(1193, 7).Value  = "Papa"
$ws.Cells.Item(1193, 8).Value  = "Rosara"
$ws.Cells.Item(1193, 9).Value  = "1a nueva(o)"
$ws.Cells.Item(1193, 10).Value = 340
$ws.Cells.Item(1193, 11).Value = 19000
$ws.Cells.Item(1193, 12).Value = 20000
$ws.Cells.Item(1193, 13).Value = 19441
$ws.Cells.Item(1193, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(1193, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(1193, 16).Value = 778
$ws.Cells.Item(1193, 17).Value = 25
$ws.Cells.Item(1193, 18).Value = "Hortaliza"
